# Update "want-to-go count" (column F) values on the "展览" and "全部类型"
# sheets to reflect the latest scrape, as described by the commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll = $wb.Worksheets.Item("全部类型")

# Each entry: row on "展览" sheet, row on "全部类型" sheet, new F value
$updates = @(
    @{ ExhibitRow = 2;  AllRow = 3;  NewValue = 11 },
    @{ ExhibitRow = 3;  AllRow = 4;  NewValue = 13162 },
    @{ ExhibitRow = 4;  AllRow = 5;  NewValue = 36 },
    @{ ExhibitRow = 5;  AllRow = 6;  NewValue = 3 },
    @{ ExhibitRow = 6;  AllRow = 7;  NewValue = 101 },
    @{ ExhibitRow = 11; AllRow = 12; NewValue = 13119 },
    @{ ExhibitRow = 12; AllRow = 13; NewValue = 317 },
    @{ ExhibitRow = 13; AllRow = 14; NewValue = 561 },
    @{ ExhibitRow = 14; AllRow = 15; NewValue = 8805 },
    @{ ExhibitRow = 15; AllRow = 16; NewValue = 7870 },
    @{ ExhibitRow = 21; AllRow = 22; NewValue = 3 },
    @{ ExhibitRow = 26; AllRow = 29; NewValue = 194 },
    @{ ExhibitRow = 27; AllRow = 30; NewValue = 66 },
    @{ ExhibitRow = 28; AllRow = 31; NewValue = 343 },
    @{ ExhibitRow = 30; AllRow = 33; NewValue = 5222 }
)

foreach ($u in $updates) {
    $wsExhibit.Range("F" + $u.ExhibitRow).Value = $u.NewValue
    $wsAll.Range("F" + $u.AllRow).Value = $u.NewValue
}
